$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 489.85
$ws.Range("I19").Value = 543.6923
$ws.Range("K19").Value = 543.6923
$ws.Range("M19").Value = -368.6923
$ws.Range("H76").Value = 4926.706
$ws.Range("I76").Value = 5272.636
$ws.Range("K76").Value = 5272.636
$ws.Range("M76").Value = -4957.636
$ws.Range("H79").Value = 4926.706
$ws.Range("I79").Value = 5272.636
$ws.Range("K79").Value = 5272.636
$ws.Range("M79").Value = -4180.636
$ws.Range("H112").Value = 1934.1666
$ws.Range("J112").Value = 2007.7028
$ws.Range("L112").Value = 6023.1084
$ws.Range("N112").Value = -8239.108400000001
$ws.Range("H116").Value = 3484.842
$ws.Range("I116").Value = 3073.1428
$ws.Range("K116").Value = 3073.1428
$ws.Range("M116").Value = 368.8571999999999
$ws.Range("H129").Value = 558263.5600000001
$ws.Range("I129").Value = 810569.6
$ws.Range("K129").Value = 2431708.8
$ws.Range("M129").Value = -2426708.8
$ws.Range("H131").Value = 5108.5
$ws.Range("I131").Value = 3774.4546
$ws.Range("J131").Value = 10000
$ws.Range("K131").Value = 11323.3638
$ws.Range("L131").Value = 30000
$ws.Range("M131").Value = -6283.363799999999
$ws.Range("N131").Value = -40080
$ws.Range("H132").Value = 13982.154
$ws.Range("I132").Value = 1232.6875
$ws.Range("K132").Value = 3698.0625
$ws.Range("M132").Value = -1168.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9263.691000000001
$ws.Range("I32").Value = 8006.403
$ws.Range("J32").Value = 19322
$ws.Range("K32").Value = 8006.403
$ws.Range("L32").Value = 19322
$ws.Range("M32").Value = -7719.403
$ws.Range("N32").Value = -19896
$ws.Range("H43").Value = 19999
$ws.Range("J43").Value = 19999
$ws.Range("L43").Value = 19999
$ws.Range("N43").Value = -20625
$ws.Range("H61").Value = 1880.8077
$ws.Range("I61").Value = 922.2778
$ws.Range("K61").Value = 922.2778
$ws.Range("M61").Value = -710.2778
$ws.Range("H102").Value = 2240
$ws.Range("I102").Value = 2360
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2360
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -738
$ws.Range("N102").Value = -5244
$ws.Range("H110").Value = 2925.7222
$ws.Range("I110").Value = 2925.7222
$ws.Range("K110").Value = 2925.7222
$ws.Range("M110").Value = -880.7222000000002
$ws.Range("H122").Value = 3645.675
$ws.Range("I122").Value = 1722.5807
$ws.Range("J122").Value = 10269.667
$ws.Range("K122").Value = 5167.742099999999
$ws.Range("L122").Value = 30809.001
$ws.Range("M122").Value = -2717.742099999999
$ws.Range("N122").Value = -35709.001
$ws.Range("H132").Value = 2681.6545
$ws.Range("I132").Value = 2210.2708
$ws.Range("K132").Value = 6630.812399999999
$ws.Range("M132").Value = -4100.812399999999
$ws.Range("H136").Value = 1880.8077
$ws.Range("I136").Value = 922.2778
$ws.Range("K136").Value = 2766.8334
$ws.Range("M136").Value = -216.8334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3265.2307
$ws.Range("I86").Value = 1735
$ws.Range("J86").Value = 5050.5
$ws.Range("K86").Value = 1735
$ws.Range("L86").Value = 5050.5
$ws.Range("M86").Value = -612
$ws.Range("N86").Value = -7296.5
$ws.Range("H89").Value = 3265.2307
$ws.Range("I89").Value = 1735
$ws.Range("J89").Value = 5050.5
$ws.Range("K89").Value = 8675
$ws.Range("L89").Value = 25252.5
$ws.Range("M89").Value = -3059
$ws.Range("N89").Value = -36484.5
$ws.Range("H94").Value = 4428.4546
$ws.Range("I94").Value = 3950.5
$ws.Range("J94").Value = 5002
$ws.Range("K94").Value = 3950.5
$ws.Range("L94").Value = 5002
$ws.Range("M94").Value = -3499.5
$ws.Range("N94").Value = -5904

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 128.41667
$ws.Range("I22").Value = 94.2
$ws.Range("J22").Value = 152.85715
$ws.Range("K22").Value = 94.2
$ws.Range("L22").Value = 152.85715
$ws.Range("M22").Value = 255.8
$ws.Range("N22").Value = -852.85715
$ws.Range("H31").Value = 2739.8647
$ws.Range("I31").Value = 1176.5238
$ws.Range("J31").Value = 4791.75
$ws.Range("K31").Value = 1176.5238
$ws.Range("L31").Value = 4791.75
$ws.Range("M31").Value = -881.5237999999999
$ws.Range("N31").Value = -5381.75
$ws.Range("H34").Value = 2739.8647
$ws.Range("I34").Value = 1176.5238
$ws.Range("J34").Value = 4791.75
$ws.Range("K34").Value = 1176.5238
$ws.Range("L34").Value = 4791.75
$ws.Range("M34").Value = -974.5237999999999
$ws.Range("N34").Value = -5195.75
$ws.Range("H134").Value = 2488.0425
$ws.Range("I134").Value = 1241.4
$ws.Range("K134").Value = 3724.2
$ws.Range("M134").Value = -1189.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 9544.412
$ws.Range("J5").Value = 13364.4
$ws.Range("L5").Value = 40093.2
$ws.Range("N5").Value = -40317.2
$ws.Range("H128").Value = 999999.7
$ws.Range("I128").Value = 999999.7
$ws.Range("K128").Value = 2999999.1
$ws.Range("M128").Value = -2995019.1
$ws.Range("H135").Value = 9544.412
$ws.Range("J135").Value = 13364.4
$ws.Range("L135").Value = 120279.6
$ws.Range("N135").Value = -125349.6
$ws.Range("H137").Value = 1758.8667
$ws.Range("I137").Value = 1563.3
$ws.Range("J137").Value = 2150
$ws.Range("K137").Value = 4689.9
$ws.Range("L137").Value = 6450
$ws.Range("M137").Value = 410.1000000000004
$ws.Range("N137").Value = -16650

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 77400.07000000001
$ws.Range("I80").Value = 114129.5
$ws.Range("J80").Value = 3941.2
$ws.Range("K80").Value = 114129.5
$ws.Range("L80").Value = 3941.2
$ws.Range("M80").Value = -113131.5
$ws.Range("N80").Value = -5937.2
$ws.Range("H83").Value = 77400.07000000001
$ws.Range("I83").Value = 114129.5
$ws.Range("J83").Value = 3941.2
$ws.Range("K83").Value = 570647.5
$ws.Range("L83").Value = 19706
$ws.Range("M83").Value = -565655.5
$ws.Range("N83").Value = -29690
$ws.Range("H97").Value = 500.85
$ws.Range("I97").Value = 434.54544
$ws.Range("J97").Value = 581.8889
$ws.Range("K97").Value = 434.54544
$ws.Range("L97").Value = 581.8889
$ws.Range("M97").Value = 61.45456000000001
$ws.Range("N97").Value = -1573.8889
$ws.Range("H102").Value = 2538.7307
$ws.Range("I102").Value = 1891.7391
$ws.Range("K102").Value = 1891.7391
$ws.Range("M102").Value = -269.7391
$ws.Range("H122").Value = 2954.6667
$ws.Range("I122").Value = 3197.2632
$ws.Range("K122").Value = 9591.7896
$ws.Range("M122").Value = -7141.7896
$ws.Range("H126").Value = 4121
$ws.Range("I126").Value = 2451.4
$ws.Range("J126").Value = 5790.6
$ws.Range("K126").Value = 7354.200000000001
$ws.Range("L126").Value = 17371.8
$ws.Range("M126").Value = -4884.200000000001
$ws.Range("N126").Value = -22311.8
$ws.Range("H132").Value = 2329.9429
$ws.Range("I132").Value = 2167.303
$ws.Range("K132").Value = 6501.909
$ws.Range("M132").Value = -3971.909

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 5555
$ws.Range("I35").Value = 5555
$ws.Range("K35").Value = 5555
$ws.Range("M35").Value = -5219
$ws.Range("H61").Value = 2894.818
$ws.Range("I61").Value = 980.4375
$ws.Range("K61").Value = 980.4375
$ws.Range("M61").Value = -778.4375
$ws.Range("H100").Value = 224562
$ws.Range("J100").Value = 3099.6667
$ws.Range("L100").Value = 3099.6667
$ws.Range("N100").Value = -4181.6667
$ws.Range("H113").Value = 2894.818
$ws.Range("I113").Value = 980.4375
$ws.Range("K113").Value = 980.4375
$ws.Range("M113").Value = 1189.5625
$ws.Range("H122").Value = 4582.4517
$ws.Range("I122").Value = 3652.348
$ws.Range("J122").Value = 7256.5
$ws.Range("K122").Value = 10957.044
$ws.Range("L122").Value = 21769.5
$ws.Range("M122").Value = -8507.044
$ws.Range("N122").Value = -26669.5
$ws.Range("H136").Value = 4177.75
$ws.Range("I136").Value = 2371.4211
$ws.Range("J136").Value = 5812.048
$ws.Range("K136").Value = 7114.263300000001
$ws.Range("L136").Value = 17436.144
$ws.Range("M136").Value = -4564.263300000001
$ws.Range("N136").Value = -22536.144

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 9899
$ws.Range("I5").Value = 9899
$ws.Range("K5").Value = 9899
$ws.Range("M5").Value = -9787
$ws.Range("H96").Value = 1960.4117
$ws.Range("I96").Value = 1777.3077
$ws.Range("J96").Value = 2555.5
$ws.Range("K96").Value = 1777.3077
$ws.Range("L96").Value = 2555.5
$ws.Range("M96").Value = -404.3077000000001
$ws.Range("N96").Value = -5301.5
$ws.Range("H136").Value = 2173.6978
$ws.Range("I136").Value = 1208.6666
$ws.Range("K136").Value = 3625.9998
$ws.Range("M136").Value = -1075.9998